# Update Pierre Data Collection
# Mark rows 72-74 in column C as "ok" to match rows 75-88, and move the
# active selection to E75.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C72").Value = "ok"
$ws.Range("C73").Value = "ok"
$ws.Range("C74").Value = "ok"

$ws.Range("E75").Select()
